$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B26").Value = "clone the git remote repository to local reposiroty"
$ws.Range("C26").Value = "git clone https://github.com/KrishnaGopalGupta/Git-Practice.git"
$ws.Range("A26").Value = "git clone <git remote url>"
$ws.Columns.Item(3).ColumnWidth = 58.6
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("B27").Select()
